# The source publishes a new weekly price record for
# "Feria Lagunitas de Puerto Montt" / Arándano (blue). It belongs right
# after the current row 7 (chronologically/logically it becomes the new
# row 8), so insert a fresh row there and push the existing rows 8-32 down
# to 9-33 (carrying their formatting, in particular the date style on
# column D, along with them).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C8").Value = "Los Lagos"
$ws.Range("D8").Value = 44910
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101001
$ws.Range("J8").Value = "Arándano (blue)"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 4000
$ws.Range("O8").Value = 4200
$ws.Range("P8").Value = 4100
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 2050
$ws.Range("T8").Value = 2
